$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 5.682999999999999
$ws.Range("D12").Value = -7.269
$ws.Range("D14").Value = -7.505
$ws.Range("B18").Value = 5.165999999999999
$ws.Range("B20").Value = 6.613000000000001
$ws.Range("D26").Value = -8.122999999999999
$ws.Range("B27").Value = 6.105
$ws.Range("D27").Value = -8.251999999999999
$ws.Range("D29").Value = -7.297999999999999
$ws.Range("B35").Value = 7.657999999999999
$ws.Range("D37").Value = -7.822
$ws.Range("D38").Value = -7.806
$ws.Range("D51").Value = -8.638000000000002
$ws.Range("D52").Value = -7.831
$ws.Range("D55").Value = -8.129000000000001
$ws.Range("B69").Value = 5.319999999999999
$ws.Range("D69").Value = -7.017999999999999
$ws.Range("D70").Value = -6.874
$ws.Range("B76").Value = 6.308
$ws.Range("B78").Value = 8.083
$ws.Range("D81").Value = -7.505000000000001
$ws.Range("B82").Value = 5.366000000000001
$ws.Range("B83").Value = 5.305
$ws.Range("D83").Value = -8.549000000000001
$ws.Range("B93").Value = 6.045
$ws.Range("D102").Value = -7.833
